$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing data row (row 2) down onto the new
# row 3 so the appended record visually matches the rest of the table
# (thin border style, default font).
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null

# Plain text columns: name / nickname / email assign normally.
$ws.Range("A3").Value = "홍길동"
$ws.Range("B3").Value = "kdw1234"
$ws.Range("D3").Value = "kdw8573@snu.ac.kr"

# Phone number must stay as literal text so the leading zero survives
# (a plain numeric-looking string assigned to a General-formatted cell
# gets auto-coerced to a number). Build it as text in a scratch cell and
# bring only the *value* back onto C3 so C3 keeps its original style/format.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01085732136"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# Last column is numeric.
$ws.Range("E3").Value = 0
